$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title "HARSHIL DALE" - bump font size from 24pt (sz/szCs=48) to 28pt
#    (sz/szCs=56). Operate on the whole paragraph Range (which also covers
#    the paragraph-mark rPr stored in pPr/rPr) so every run (and the mark)
#    picks up both the ASCII size (Font.Size -> w:sz) and the complex-script
#    size (Font.SizeBi -> w:szCs).
# ---------------------------------------------------------------------------
$title = $d.Paragraphs(1).Range
if ($title.Text -notmatch "HARSHIL") {
    throw "Unexpected first paragraph text: $($title.Text)"
}
$title.Font.Size = 28
$title.Font.SizeBi = 28

# ---------------------------------------------------------------------------
# 2) Merge the "Delivered high-impact lectures ... as a Lecturer" run with
#    the trailing "." run into a single run (same text, same formatting).
#    We locate the sentence (without the final period) and the following
#    period character, drop the period, then re-append it right onto the
#    end of the first run so it becomes part of that run instead of living
#    in a run of its own.
# ---------------------------------------------------------------------------
$sentence = $d.Content
$found = $sentence.Find.Execute("Delivered high-impact lectures and guided data-focused student projects as a Lecturer", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Lecturer sentence"
}

$sentenceEnd = $sentence.End
$period = $d.Range($sentenceEnd, $sentenceEnd + 1)
if ($period.Text -ne ".") {
    throw "Expected a trailing period after the Lecturer sentence, found: $($period.Text)"
}
$period.Delete()
$sentence.InsertAfter(".")

Write-Output "done"
